$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Replacement Screen and Solder"
$ws.Range("B12").Value = 20.44
$ws.Range("C12").Value = "NA"

$ws.ListObjects.Item(1).Resize($ws.Range("A1:C13"))

$ws.Range("C18").Select()
